$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "1.0000", "4.320", "0.8410") are preserved verbatim as text,
# matching the source inlineStr cells instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.375.59"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.880.84"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "244.78"
$ws.Range("E5").Value = "  +4.55%  "
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.4774"
$ws.Range("E7").Value = "  +1.89%  "
$ws.Range("D8").Value = "0.2882"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").Value = "0.06527"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "21.39"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").Value = "0.07763"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").Value = "96.74"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.880.60"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").Value = "0.7363"
$ws.Range("E14").Value = "  +7.00%  "
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "274.99"
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("D17").Value = "30.364.42"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "0.000007558"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "2.137.60"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "6.178"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("E25").Value = "  -2.02%  "
$ws.Range("D26").Value = "163.59"
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("D27").Value = "18.94"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "1.964"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").Value = "1.372"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").Value = "0.09989"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "1.511"
$ws.Range("D32").Value = "4.320"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "4.083"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "0.04747"
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "0.6964"
$ws.Range("D37").Value = "2.717"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "2.749"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").Value = "6.271"
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("D41").Value = "69.49"
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").Value = "0.8410"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").Value = "0.9993"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "101.81"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.295"
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "7.094"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "35.21"
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").Value = "911.46"
$ws.Range("E50").Value = "  -5.56%  "
$ws.Range("D51").Value = "0.05594"
$ws.Range("E51").Value = "  -0.72%  "
